$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.919.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.218.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '292.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.469'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0781'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '50.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.19%  '
$ws.Range("E13").Value = '  +2.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.556.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.79'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.244.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.732'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.863.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0887'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0711'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.47%  '
$ws.Range("E37").Value = '  -1.64%  '
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0986'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("E40").Value = '  +1.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.35'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.110.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0270'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.427.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
